$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.408.90"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.722.22"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.86"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4860"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2584"
$ws.Range("E8").Value = "  -3.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06190"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.718.23"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06975"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.45"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.540"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5966"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.23"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.413.23"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007224"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.30"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.949.41"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.456"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.494"
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.05"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.400"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.78"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.724"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.945"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07973"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.663"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04495"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9994"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.600"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9975"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6263"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9340"
$ws.Range("E38").Value = "  +4.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.953"
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9999"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01472"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.70"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.833"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05359"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.750"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.19"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.227"
$ws.Range("E51").Value = "  -1.96%  "
